$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = 42607.891712962963

$ws.Range("B4").Value = -38
$ws.Range("C4").Value = 36
$ws.Range("D4").Value = 62
$ws.Range("E4").Value = 0
$ws.Range("F4").Value = 100
$ws.Range("G4").Value = 19530
$ws.Range("H4").Value = 16393
$ws.Range("I4").Value = 810
$ws.Range("J4").Value = 125
$ws.Range("K4").Value = 214
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = 10
$ws.Range("N4").Value = "Named"
